# Apply the price-list update described by the commit diff.
# Rows 3-12 (the "iPhone 14 Pro Max ..." section) get new D-column prices,
# their E-column ("Price with 4% ?") formulas recompute automatically,
# and the country labels in column A become bold for that block.
# A couple of rows also get a bespoke (non-4%) markup formula and one
# country value is corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: bold the country/model labels for rows 3 through 12 ---
3..12 | ForEach-Object {
    $ws.Cells.Item($_, 1).Font.Bold = $true
}

# --- Row 2: price bump + special 1% markup formula (no longer shared) ---
$ws.Range("D2").Value = 150500
$ws.Range("E2").Formula = "=(D2 * 0.01) + D2"

# --- Row 3: price bump, standard 4% markup ---
$ws.Range("D3").Value = 125000

# --- Row 5: price bump ---
$ws.Range("D5").Value = 137000

# --- Row 6: price bump ---
$ws.Range("D6").Value = 135750

# --- Row 7: price bump ---
$ws.Range("D7").Value = 105700

# --- Row 8: price bump ---
$ws.Range("D8").Value = 105800

# --- Row 9: price bump ---
$ws.Range("D9").Value = 95500

# --- Row 10: price bump ---
$ws.Range("D10").Value = 95500

# --- Row 11: country corrected to Saudi Arabia, price bump, special 5% markup ---
$ws.Range("C11").Value = "Saudi Arabia 🇸🇦"
$ws.Range("D11").Value = 99000
$ws.Range("E11").Formula = "=(D11 * 0.05) + D11"

# --- Row 12: price bump, special (explicit) 4% markup formula ---
$ws.Range("D12").Value = 94900
$ws.Range("E12").Formula = "=(D12 * 0.04) + D12"
